# Add data for 2022-02-22: update the "through" date label and bump the
# underlying carjacking counts for the neighborhoods affected by the new
# day's records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet and update the matching header cell/shared string.
$ws.Name = "Through 2022-02-14"
$ws.Range("B1").Value = "February 2022 (through February 14)"

# Row 2 - Englewood
$ws.Range("F2").Value = 2
$ws.Range("J2").Value = 1

# Row 3 - Austin
$ws.Range("J3").Value = 5
$ws.Range("L3").Value = 5

# Row 4 - New City
$ws.Range("N4").Value = 2

# Row 6 - South Shore
$ws.Range("N6").Value = 1

# Row 7 - Auburn Gresham
$ws.Range("D7").Value = 2

# Row 13 - Bridgeport
$ws.Range("B13").Value = 2

# Row 15 - Garfield Park
$ws.Range("B15").Value = 5

# Row 18 - Humboldt Park
$ws.Range("D18").Value = 1

# Row 19 - Logan Square
$ws.Range("N19").Value = 1

# Row 22 - Chicago Lawn
$ws.Range("N22").Value = 4

# Row 38 - Wicker Park
$ws.Range("B38").Value = 1

# Row 48 - Lake View
$ws.Range("B48").Value = 3

# Row 57 - Woodlawn
$ws.Range("D57").Value = 1

# Row 60 - Chinatown
$ws.Range("B60").Value = 4
